$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'23.359.29"
$ws.Range("E2").Value2 = "  -0.86%  "

$ws.Range("D3").Value2 = "'1.625.59"
$ws.Range("E3").Value2 = "  -0.83%  "

$ws.Range("D4").Value2 = "'0.9993"
$ws.Range("E4").Value2 = "  +0.23%  "

$ws.Range("D5").Value2 = "'1.000"
$ws.Range("E5").Value2 = "  +0.25%  "

$ws.Range("D6").Value2 = "'304.08"
$ws.Range("E6").Value2 = "  -1.31%  "

$ws.Range("D7").Value2 = "'0.3781"

$ws.Range("D8").Value2 = "'52.06"
$ws.Range("E8").Value2 = "  -2.16%  "

$ws.Range("D9").Value2 = "'0.3613"
$ws.Range("E9").Value2 = "  -1.71%  "

$ws.Range("B10").Value2 = "Dogecoin"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value2 = "'0.08094"
$ws.Range("E10").Value2 = "  -1.29%  "

$ws.Range("B11").Value2 = "Polygon"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value2 = "'1.226"
$ws.Range("E11").Value2 = "  -4.74%  "

$ws.Range("D12").Value2 = "'1.001"
$ws.Range("E12").Value2 = "  +0.38%  "

$ws.Range("D13").Value2 = "'22.66"
$ws.Range("E13").Value2 = "  -2.93%  "

$ws.Range("D14").Value2 = "'6.553"
$ws.Range("E14").Value2 = "  -1.95%  "

$ws.Range("D15").Value2 = "'0.00001246"
$ws.Range("E15").Value2 = "  -3.19%  "

$ws.Range("D16").Value2 = "'7.216"
$ws.Range("E16").Value2 = "  -3.70%  "

$ws.Range("D17").Value2 = "'1.621.16"
$ws.Range("E17").Value2 = "  -0.62%  "

$ws.Range("D18").Value2 = "'93.53"
$ws.Range("E18").Value2 = "  -1.59%  "

$ws.Range("D19").Value2 = "'0.06908"
$ws.Range("E19").Value2 = "  -0.73%  "

$ws.Range("D20").Value2 = "'17.87"
$ws.Range("E20").Value2 = "  -3.39%  "

$ws.Range("D21").Value2 = "'1.001"
$ws.Range("E21").Value2 = "  +0.30%  "

$ws.Range("D22").Value2 = "'6.417"
$ws.Range("E22").Value2 = "  -2.95%  "

$ws.Range("D23").Value2 = "'23.366.09"
$ws.Range("E23").Value2 = "  -0.80%  "

$ws.Range("D24").Value2 = "'12.68"
$ws.Range("E24").Value2 = "  -2.66%  "

$ws.Range("D25").Value2 = "'3.208"
$ws.Range("E25").Value2 = "  +1.98%  "

$ws.Range("D26").Value2 = "'2.431"
$ws.Range("E26").Value2 = "  +0.38%  "

$ws.Range("E27").Value2 = "  -1.88%  "

$ws.Range("D28").Value2 = "'149.63"
$ws.Range("E28").Value2 = "  -1.29%  "

$ws.Range("D29").Value2 = "'5.285"
$ws.Range("E29").Value2 = "  -0.70%  "

$ws.Range("D30").Value2 = "'134.56"

$ws.Range("D31").Value2 = "'2.295"
$ws.Range("E31").Value2 = "  -5.58%  "

$ws.Range("D32").Value2 = "'1.802.15"
$ws.Range("E32").Value2 = "  -0.50%  "

$ws.Range("D33").Value2 = "'6.768"
$ws.Range("E33").Value2 = "  -1.50%  "

$ws.Range("D34").Value2 = "'11.00"
$ws.Range("E34").Value2 = "  +5.18%  "

$ws.Range("D35").Value2 = "'0.9490"
$ws.Range("E35").Value2 = "  -3.02%  "

$ws.Range("D36").Value2 = "'0.02770"
$ws.Range("E36").Value2 = "  -2.13%  "

$ws.Range("D37").Value2 = "'0.2513"
$ws.Range("E37").Value2 = "  -1.22%  "

$ws.Range("D38").Value2 = "'0.08816"
$ws.Range("E38").Value2 = "  -0.51%  "

$ws.Range("B39").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value2 = "'6.057"
$ws.Range("E39").Value2 = "  -2.98%  "

$ws.Range("B40").Value2 = "Hedera"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value2 = "'0.07111"
$ws.Range("E40").Value2 = "  -5.27%  "

$ws.Range("D41").Value2 = "'1.359"
$ws.Range("E41").Value2 = "  -2.99%  "

$ws.Range("D42").Value2 = "'0.7053"
$ws.Range("E42").Value2 = "  -1.91%  "

$ws.Range("D43").Value2 = "'16.14"
$ws.Range("E43").Value2 = "  -0.46%  "

$ws.Range("D44").Value2 = "'12.27"
$ws.Range("E44").Value2 = "  -3.77%  "

$ws.Range("D45").Value2 = "'0.6441"
$ws.Range("E45").Value2 = "  -3.35%  "

$ws.Range("B46").Value2 = "NEARProtocol"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value2 = "'2.313"
$ws.Range("E46").Value2 = "  -2.51%  "

$ws.Range("B47").Value2 = "Frax"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value2 = "'0.9996"
$ws.Range("E47").Value2 = "  +0.34%  "

$ws.Range("D48").Value2 = "'3.991"
$ws.Range("E48").Value2 = "  -1.22%  "

$ws.Range("D49").Value2 = "'0.07977"
$ws.Range("E49").Value2 = "  -0.99%  "

$ws.Range("D50").Value2 = "'1.195"
$ws.Range("E50").Value2 = "  -1.84%  "

$ws.Range("D51").Value2 = "'125.67"
$ws.Range("E51").Value2 = "  -4.86%  "
